$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Last modification" date in the Goal outline / Level doc info table
#    changes from 5/12/2017 to 29/01/2018.
#    Target the specific table cell directly (located by its row label)
#    so the unrelated "05/12/2017" entry in the Revision History table
#    (and the "29/01/2018" entry already present there) are left
#    untouched.
# ---------------------------------------------------------------------
foreach ($t in $d.Tables) {
    for ($r = 1; $r -le $t.Rows.Count; $r++) {
        $label = $t.Cell($r, 1).Range.Text
        if ($label -like "Last modification*") {
            $valueCell = $t.Cell($r, 2)
            if ($valueCell.Range.Text -like "5/12/2017*") {
                $valueCell.Range.Text = "29/01/2018"
            }
        }
    }
}

# ---------------------------------------------------------------------
# 2) Backup section rewording:
#    ", the master backup is kept by Sara Obici."
#    -> ", the master backup, for this reason, is kept by github itself."
# ---------------------------------------------------------------------
$oldSentence = ", the master backup is kept by Sara Obici."
$newSentence = ", the master backup, for this reason, is kept by github itself."
$d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, $true, 1, $false, $newSentence, 2) | Out-Null

# ---------------------------------------------------------------------
# 3) The "_GoBack" bookmark moves from its old spot (before the
#    "Development, which contains..." bullet) to right before the
#    second "github" occurrence we just inserted above.
#    Re-adding a bookmark named "_GoBack" relocates the (singleton)
#    bookmark, removing it from its previous position automatically.
# ---------------------------------------------------------------------
$targetRange = $d.Content
$targetRange.Find.Execute("github itself", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$bookmarkRange = $d.Range($targetRange.Start, $targetRange.Start)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# ---------------------------------------------------------------------
# 4) Remove the whole bullet paragraph:
#    "Last master backup: on an external memory (Toshiba) 9/1/2018"
# ---------------------------------------------------------------------
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Last master backup: on an external memory (Toshiba) 9/1/2018*") {
        $p.Range.Delete()
        break
    }
}
